# Insert a new data row at row 12 (shifting existing rows 12..93 down to 13..94)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12; existing rows shift down automatically.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44635
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112030
$ws.Cells.Item(12, 7).Value = "Poroto granado"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 20000
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 20000
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(12, 16).Value = 800
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# Match the style used by the other date cells in column D (numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
